# Generate Report for Handback
#
# The localization-status report is refreshed once a handback has been
# received and verified as in sync with the en-US source: the "Status"
# columns move from "Ready for handoff" to "Handed back: in sync with
# en-US", the per-locale "Latest Handback DateTime" stamps are bumped to
# the time the handback was recorded, and the stale "Error Detail" note
# (which complained the handback file wasn't the latest revision) is
# cleared now that everything is current. Column widths are re-fit to the
# new cell contents.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview.Range("E2").Value = $status
$overview.Range("F2").Value = $status

# zh-cn detail sheet: status, refreshed handback datetime, clear error detail
$zhcn.Range("C2").Value = $status
$zhcn.Range("K2").Value = "2016-08-25 18:51:29"
$zhcn.Range("P2").Value = ""

# de-de detail sheet: status, refreshed handback datetime, clear error detail
$dede.Range("C2").Value = $status
$dede.Range("K2").Value = "2016-08-25 18:51:36"
$dede.Range("P2").Value = ""

# Re-fit the columns whose contents changed length
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
